# Apply updated Cosinor analysis results (rows 2 and 3) per re-run of
# CircaDB / CircadiPy simulation analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 24.52000000000039
$ws.Range("G2").Value = [double]"1.909868929672598e-09"
$ws.Range("H2").Value = [double]"5.598230269377345e-09"
$ws.Range("K2").Value = 34.42829519000485
$ws.Range("L2").Value = "[22.134984907096815, 46.72160547291288]"
$ws.Range("M2").Value = [double]"1.083992438033476e-07"
$ws.Range("N2").Value = [double]"1.083992438033476e-07"
$ws.Range("O2").Value = 1.415131825941348
$ws.Range("P2").Value = "[1.0126054398958093, 1.8176582119868865]"
$ws.Range("Q2").Value = [double]"6.264255780763506e-11"
$ws.Range("R2").Value = [double]"6.264255780763506e-11"
$ws.Range("S2").Value = 57.34599405497261
$ws.Range("T2").Value = "[50.12737150612605, 64.56461660381916]"
$ws.Range("W2").Value = 18.99747747747778
$ws.Range("X2").Value = 17.42662662662691
$ws.Range("Y2").Value = 20.56832832832866

# --- Row 3 ---
$ws.Range("E3").Value = 25.8200000000006
$ws.Range("H3").Value = [double]"5.255493607693049e-16"
$ws.Range("K3").Value = 56.50628849184892
$ws.Range("L3").Value = "[43.97825884729704, 69.0343181364008]"
$ws.Range("M3").Value = [double]"2.220446049250313e-16"
$ws.Range("N3").Value = [double]"4.440892098500626e-16"
$ws.Range("O3").Value = 2.572395185822273
$ws.Range("P3").Value = "[2.3459740936716575, 2.7988162779728882]"
$ws.Range("S3").Value = 67.64898558371334
$ws.Range("T3").Value = "[61.05051838888696, 74.24745277853972]"
$ws.Range("W3").Value = 15.2490490490494
$ws.Range("X3").Value = 14.31859859859893
$ws.Range("Y3").Value = 16.17949949949987
